# Generate Report for Handoff
# The e2e\324a4567-38de-4d2d-975e-9d6a52a3674c.md file has been handed off for
# translation of a00a7228-422b-48f1-b114-67c1f80c027f.md, so the two rows
# swap places (324a4567 now sorts above a00a7228) and the a00a7228 row
# picks up the new "Ready for handoff" status together with fresh handoff
# timestamps.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview" (A1:G9)
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

# Row 8 becomes what row 9 used to be (324a4567...)
$overview.Range("A8").Value = "324a4567-38de-4d2d-975e-9d6a52a3674c.md"
$overview.Range("B8").Value = "e2e\324a4567-38de-4d2d-975e-9d6a52a3674c.md"
$overview.Range("E8").Value = "Ready for handoff"
$overview.Range("F8").Value = "Ready for handoff"
$overview.Range("G8").Value = "2016-10-19 16:39:11"

# Row 9 becomes what row 8 used to be (a00a7228...), with an updated handoff
$overview.Range("A9").Value = "a00a7228-422b-48f1-b114-67c1f80c027f.md"
$overview.Range("B9").Value = "e2e\a00a7228-422b-48f1-b114-67c1f80c027f.md"
$overview.Range("E9").Value = "Ready for handoff"
$overview.Range("F9").Value = "Ready for handoff"
$overview.Range("G9").Value = "2016-10-19 16:49:16"

# Rebuild the hyperlinks for column B (deleting one hyperlink clears the
# whole collection in this engine, so clear once and re-add every link).
$overview.Range("A1").Hyperlinks.Delete()
$overview.Hyperlinks.Add($overview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3c0e21910103a15be3449acf017a5abab43b3ec2/e2e/44d6cb8a-8333-4c77-b709-5f468dd28896.md", [Type]::Missing, [Type]::Missing, "e2e\44d6cb8a-8333-4c77-b709-5f468dd28896.md")
$overview.Hyperlinks.Add($overview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5864a086a2d4623c58f3ad011ba79cba79b30fd0/e2e/1692af14-de8b-457c-b89b-371ae2d85f56.md", [Type]::Missing, [Type]::Missing, "e2e\1692af14-de8b-457c-b89b-371ae2d85f56.md")
$overview.Hyperlinks.Add($overview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/373fe962133c8fcac57e6678edee7deadbecc1cc/e2e/3304663f-245b-4523-a4de-244871a32b5b.md", [Type]::Missing, [Type]::Missing, "e2e\3304663f-245b-4523-a4de-244871a32b5b.md")
$overview.Hyperlinks.Add($overview.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b8dca5e31295c1e8f1d2eff6a350d61283c92c6f/e2e/4eed8814-8aed-4dd0-ae74-9c5148eb643d.md", [Type]::Missing, [Type]::Missing, "e2e\4eed8814-8aed-4dd0-ae74-9c5148eb643d.md")
$overview.Hyperlinks.Add($overview.Range("B6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/adbd0af2d4e2f8b95939cad31de737fc83aadf2c/e2e/5c4c6826-7756-4723-a923-e65d0f2de573.md", [Type]::Missing, [Type]::Missing, "e2e\5c4c6826-7756-4723-a923-e65d0f2de573.md")
$overview.Hyperlinks.Add($overview.Range("B7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b8dca5e31295c1e8f1d2eff6a350d61283c92c6f/e2e/8f455da8-3134-4036-ac48-b5d5292b4f05.md", [Type]::Missing, [Type]::Missing, "e2e\8f455da8-3134-4036-ac48-b5d5292b4f05.md")
$overview.Hyperlinks.Add($overview.Range("B8"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5accfe475a80f6da58cc05cce09d320d0fa4319b/e2e/324a4567-38de-4d2d-975e-9d6a52a3674c.md", [Type]::Missing, [Type]::Missing, "e2e\324a4567-38de-4d2d-975e-9d6a52a3674c.md")
$overview.Hyperlinks.Add($overview.Range("B9"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f8ad85cd858171fa98d486127ccb581304fd77e5/e2e/a00a7228-422b-48f1-b114-67c1f80c027f.md", [Type]::Missing, [Type]::Missing, "e2e\a00a7228-422b-48f1-b114-67c1f80c027f.md")

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (A1:P9)
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("A8").Value = "324a4567-38de-4d2d-975e-9d6a52a3674c.md"
$zhcn.Range("C8").Value = "Ready for handoff"
$zhcn.Range("G8").Value = "324a4567-38de-4d2d-975e-9d6a52a3674c.5e02424a11a8004174b34e7fcb9bc4a1b236430b.zh-cn.xlf"
$zhcn.Range("H8").Value = "2016-10-19 16:38:59"

$zhcn.Range("A9").Value = "a00a7228-422b-48f1-b114-67c1f80c027f.md"
$zhcn.Range("C9").Value = "Ready for handoff"
$zhcn.Range("G9").Value = "a00a7228-422b-48f1-b114-67c1f80c027f.73f25297e0291b637cf1e62f058f73fb1aaa5130.zh-cn.xlf"
$zhcn.Range("H9").Value = "2016-10-19 16:49:05"

$zhcn.Range("A1").Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3c0e21910103a15be3449acf017a5abab43b3ec2/e2e/44d6cb8a-8333-4c77-b709-5f468dd28896.md", [Type]::Missing, [Type]::Missing, "44d6cb8a-8333-4c77-b709-5f468dd28896.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/b101b20b54a2959c69e5410ad5e73e5dfd7b07b5/e2e/44d6cb8a-8333-4c77-b709-5f468dd28896.md", [Type]::Missing, [Type]::Missing, "44d6cb8a-8333-4c77-b709-5f468dd28896.md")
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5864a086a2d4623c58f3ad011ba79cba79b30fd0/e2e/1692af14-de8b-457c-b89b-371ae2d85f56.md", [Type]::Missing, [Type]::Missing, "1692af14-de8b-457c-b89b-371ae2d85f56.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/79d561490cb44d8c85ae944bd4fb7edb435e0cd8/e2e/1692af14-de8b-457c-b89b-371ae2d85f56.md", [Type]::Missing, [Type]::Missing, "1692af14-de8b-457c-b89b-371ae2d85f56.md")
$zhcn.Hyperlinks.Add($zhcn.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/373fe962133c8fcac57e6678edee7deadbecc1cc/e2e/3304663f-245b-4523-a4de-244871a32b5b.md", [Type]::Missing, [Type]::Missing, "3304663f-245b-4523-a4de-244871a32b5b.md")
$zhcn.Hyperlinks.Add($zhcn.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b8dca5e31295c1e8f1d2eff6a350d61283c92c6f/e2e/4eed8814-8aed-4dd0-ae74-9c5148eb643d.md", [Type]::Missing, [Type]::Missing, "4eed8814-8aed-4dd0-ae74-9c5148eb643d.md")
$zhcn.Hyperlinks.Add($zhcn.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/adbd0af2d4e2f8b95939cad31de737fc83aadf2c/e2e/5c4c6826-7756-4723-a923-e65d0f2de573.md", [Type]::Missing, [Type]::Missing, "5c4c6826-7756-4723-a923-e65d0f2de573.md")
$zhcn.Hyperlinks.Add($zhcn.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b8dca5e31295c1e8f1d2eff6a350d61283c92c6f/e2e/8f455da8-3134-4036-ac48-b5d5292b4f05.md", [Type]::Missing, [Type]::Missing, "8f455da8-3134-4036-ac48-b5d5292b4f05.md")
$zhcn.Hyperlinks.Add($zhcn.Range("A8"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5accfe475a80f6da58cc05cce09d320d0fa4319b/e2e/324a4567-38de-4d2d-975e-9d6a52a3674c.md", [Type]::Missing, [Type]::Missing, "324a4567-38de-4d2d-975e-9d6a52a3674c.md")
$zhcn.Hyperlinks.Add($zhcn.Range("A9"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f8ad85cd858171fa98d486127ccb581304fd77e5/e2e/a00a7228-422b-48f1-b114-67c1f80c027f.md", [Type]::Missing, [Type]::Missing, "a00a7228-422b-48f1-b114-67c1f80c027f.md")

# ---------------------------------------------------------------------------
# Sheet "de-de" (A1:P9)
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("A8").Value = "324a4567-38de-4d2d-975e-9d6a52a3674c.md"
$dede.Range("C8").Value = "Ready for handoff"
$dede.Range("G8").Value = "324a4567-38de-4d2d-975e-9d6a52a3674c.5e02424a11a8004174b34e7fcb9bc4a1b236430b.de-de.xlf"
$dede.Range("H8").Value = "2016-10-19 16:39:11"

$dede.Range("A9").Value = "a00a7228-422b-48f1-b114-67c1f80c027f.md"
$dede.Range("C9").Value = "Ready for handoff"
$dede.Range("G9").Value = "a00a7228-422b-48f1-b114-67c1f80c027f.73f25297e0291b637cf1e62f058f73fb1aaa5130.de-de.xlf"
$dede.Range("H9").Value = "2016-10-19 16:49:16"

$dede.Range("A1").Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3c0e21910103a15be3449acf017a5abab43b3ec2/e2e/44d6cb8a-8333-4c77-b709-5f468dd28896.md", [Type]::Missing, [Type]::Missing, "44d6cb8a-8333-4c77-b709-5f468dd28896.md")
$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/d3a2a4b23352f9869debe874d8c3355259988ac7/e2e/44d6cb8a-8333-4c77-b709-5f468dd28896.md", [Type]::Missing, [Type]::Missing, "44d6cb8a-8333-4c77-b709-5f468dd28896.md")
$dede.Hyperlinks.Add($dede.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5864a086a2d4623c58f3ad011ba79cba79b30fd0/e2e/1692af14-de8b-457c-b89b-371ae2d85f56.md", [Type]::Missing, [Type]::Missing, "1692af14-de8b-457c-b89b-371ae2d85f56.md")
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/d8833d8c62dc721424f0d21a847c9c53096451d7/e2e/1692af14-de8b-457c-b89b-371ae2d85f56.md", [Type]::Missing, [Type]::Missing, "1692af14-de8b-457c-b89b-371ae2d85f56.md")
$dede.Hyperlinks.Add($dede.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/373fe962133c8fcac57e6678edee7deadbecc1cc/e2e/3304663f-245b-4523-a4de-244871a32b5b.md", [Type]::Missing, [Type]::Missing, "3304663f-245b-4523-a4de-244871a32b5b.md")
$dede.Hyperlinks.Add($dede.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b8dca5e31295c1e8f1d2eff6a350d61283c92c6f/e2e/4eed8814-8aed-4dd0-ae74-9c5148eb643d.md", [Type]::Missing, [Type]::Missing, "4eed8814-8aed-4dd0-ae74-9c5148eb643d.md")
$dede.Hyperlinks.Add($dede.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/adbd0af2d4e2f8b95939cad31de737fc83aadf2c/e2e/5c4c6826-7756-4723-a923-e65d0f2de573.md", [Type]::Missing, [Type]::Missing, "5c4c6826-7756-4723-a923-e65d0f2de573.md")
$dede.Hyperlinks.Add($dede.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b8dca5e31295c1e8f1d2eff6a350d61283c92c6f/e2e/8f455da8-3134-4036-ac48-b5d5292b4f05.md", [Type]::Missing, [Type]::Missing, "8f455da8-3134-4036-ac48-b5d5292b4f05.md")
$dede.Hyperlinks.Add($dede.Range("A8"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5accfe475a80f6da58cc05cce09d320d0fa4319b/e2e/324a4567-38de-4d2d-975e-9d6a52a3674c.md", [Type]::Missing, [Type]::Missing, "324a4567-38de-4d2d-975e-9d6a52a3674c.md")
$dede.Hyperlinks.Add($dede.Range("A9"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f8ad85cd858171fa98d486127ccb581304fd77e5/e2e/a00a7228-422b-48f1-b114-67c1f80c027f.md", [Type]::Missing, [Type]::Missing, "a00a7228-422b-48f1-b114-67c1f80c027f.md")
